$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 357
$ws.Range("I33").Value = 189.58333
$ws.Range("J33").Value = 1026.6666
$ws.Range("K33").Value = 189.58333
$ws.Range("L33").Value = 1026.6666
$ws.Range("M33").Value = 39.41667000000001
$ws.Range("N33").Value = -1484.6666
$ws.Range("H76").Value = 11169.23
$ws.Range("I76").Value = 12381.818
$ws.Range("J76").Value = 4500
$ws.Range("K76").Value = 12381.818
$ws.Range("L76").Value = 4500
$ws.Range("M76").Value = -12066.818
$ws.Range("N76").Value = -5130
$ws.Range("H79").Value = 11169.23
$ws.Range("I79").Value = 12381.818
$ws.Range("J79").Value = 4500
$ws.Range("K79").Value = 12381.818
$ws.Range("L79").Value = 4500
$ws.Range("M79").Value = -11289.818
$ws.Range("N79").Value = -6684
$ws.Range("H95").Value = 55500
$ws.Range("J95").Value = 55500
$ws.Range("L95").Value = 55500
$ws.Range("N95").Value = -60992
$ws.Range("H108").Value = 38226
$ws.Range("J108").Value = 38226
$ws.Range("L108").Value = 38226
$ws.Range("N108").Value = -45906
$ws.Range("H112").Value = 2181.2
$ws.Range("I112").Value = 666
$ws.Range("J112").Value = 2304.054
$ws.Range("K112").Value = 1998
$ws.Range("L112").Value = 6912.162
$ws.Range("M112").Value = -890
$ws.Range("N112").Value = -9128.162
$ws.Range("H123").Value = 39371.934
$ws.Range("J123").Value = 39371.934
$ws.Range("L123").Value = 39371.934
$ws.Range("N123").Value = -49171.934
$ws.Range("H124").Value = 39847.145
$ws.Range("J124").Value = 39847.145
$ws.Range("L124").Value = 39847.145
$ws.Range("N124").Value = -49667.145
$ws.Range("H134").Value = 55333.332
$ws.Range("J134").Value = 55333.332
$ws.Range("L134").Value = 55333.332
$ws.Range("N134").Value = -65473.332
$ws.Range("H141").Value = 3896
$ws.Range("I141").Value = 3650
$ws.Range("K141").Value = 10950
$ws.Range("M141").Value = -5770

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3113.1304
$ws.Range("I2").Value = 1176.9231
$ws.Range("J2").Value = 5630.2
$ws.Range("K2").Value = 1176.9231
$ws.Range("L2").Value = 5630.2
$ws.Range("M2").Value = -1063.9231
$ws.Range("N2").Value = -5856.2
$ws.Range("H116").Value = 3113.1304
$ws.Range("I116").Value = 1176.9231
$ws.Range("J116").Value = 5630.2
$ws.Range("K116").Value = 1176.9231
$ws.Range("L116").Value = 5630.2
$ws.Range("M116").Value = 1117.0769
$ws.Range("N116").Value = -10218.2
$ws.Range("H123").Value = 662490
$ws.Range("J123").Value = 662490
$ws.Range("L123").Value = 662490
$ws.Range("N123").Value = -672290
$ws.Range("H128").Value = 52370
$ws.Range("J128").Value = 52370
$ws.Range("L128").Value = 52370
$ws.Range("N128").Value = -62330
$ws.Range("H129").Value = 39079.75
$ws.Range("J129").Value = 39079.75
$ws.Range("L129").Value = 39079.75
$ws.Range("N129").Value = -49079.75
$ws.Range("H130").Value = 42149.6
$ws.Range("J130").Value = 42149.6
$ws.Range("L130").Value = 42149.6
$ws.Range("N130").Value = -52189.6
$ws.Range("H131").Value = 38677.125
$ws.Range("J131").Value = 38677.125
$ws.Range("L131").Value = 38677.125
$ws.Range("N131").Value = -48757.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3113.1304
$ws.Range("I3").Value = 1176.9231
$ws.Range("J3").Value = 5630.2
$ws.Range("K3").Value = 1176.9231
$ws.Range("L3").Value = 5630.2
$ws.Range("M3").Value = -1062.9231
$ws.Range("N3").Value = -5858.2
$ws.Range("H122").Value = 49475
$ws.Range("J122").Value = 49475
$ws.Range("L122").Value = 49475
$ws.Range("N122").Value = -59275
$ws.Range("H124").Value = 42953.332
$ws.Range("J124").Value = 42953.332
$ws.Range("L124").Value = 42953.332
$ws.Range("N124").Value = -52773.332
$ws.Range("H129").Value = 49206.332
$ws.Range("J129").Value = 49206.332
$ws.Range("L129").Value = 49206.332
$ws.Range("N129").Value = -59206.332
$ws.Range("H130").Value = 39070.91
$ws.Range("J130").Value = 39070.91
$ws.Range("L130").Value = 39070.91
$ws.Range("N130").Value = -49110.91

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2072.1333
$ws.Range("I31").Value = 1593.8948
$ws.Range("J31").Value = 2898.182
$ws.Range("K31").Value = 1593.8948
$ws.Range("L31").Value = 2898.182
$ws.Range("M31").Value = -1298.8948
$ws.Range("N31").Value = -3488.182
$ws.Range("H34").Value = 2072.1333
$ws.Range("I34").Value = 1593.8948
$ws.Range("J34").Value = 2898.182
$ws.Range("K34").Value = 1593.8948
$ws.Range("L34").Value = 2898.182
$ws.Range("M34").Value = -1391.8948
$ws.Range("N34").Value = -3302.182
$ws.Range("H58").Value = 1425.7576
$ws.Range("I58").Value = 1345.3636
$ws.Range("J58").Value = 1586.5454
$ws.Range("K58").Value = 1345.3636
$ws.Range("L58").Value = 1586.5454
$ws.Range("M58").Value = -1142.3636
$ws.Range("N58").Value = -1992.5454
$ws.Range("H74").Value = 34382
$ws.Range("J74").Value = 34382
$ws.Range("L74").Value = 34382
$ws.Range("N74").Value = -36130
$ws.Range("H77").Value = 34382
$ws.Range("J77").Value = 34382
$ws.Range("L77").Value = 103146
$ws.Range("N77").Value = -111882
$ws.Range("H99").Value = 1436.3077
$ws.Range("I99").Value = 1493.6666
$ws.Range("J99").Value = 1387.1428
$ws.Range("K99").Value = 1493.6666
$ws.Range("L99").Value = 1387.1428
$ws.Range("M99").Value = 4.333399999999983
$ws.Range("N99").Value = -4383.1428
$ws.Range("H126").Value = 1436.3077
$ws.Range("I126").Value = 1493.6666
$ws.Range("J126").Value = 1387.1428
$ws.Range("K126").Value = 4480.9998
$ws.Range("L126").Value = 4161.428400000001
$ws.Range("M126").Value = -2010.9998
$ws.Range("N126").Value = -9101.428400000001
$ws.Range("H127").Value = 50545
$ws.Range("J127").Value = 50545
$ws.Range("L127").Value = 50545
$ws.Range("N127").Value = -60465
$ws.Range("H129").Value = 49989.5
$ws.Range("J129").Value = 49989.5
$ws.Range("L129").Value = 49989.5
$ws.Range("N129").Value = -59989.5
$ws.Range("H130").Value = 53614
$ws.Range("J130").Value = 53614
$ws.Range("L130").Value = 53614
$ws.Range("N130").Value = -63654
$ws.Range("H136").Value = 1425.7576
$ws.Range("I136").Value = 1345.3636
$ws.Range("J136").Value = 1586.5454
$ws.Range("K136").Value = 4036.0908
$ws.Range("L136").Value = 4759.6362
$ws.Range("M136").Value = -1486.0908
$ws.Range("N136").Value = -9859.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H104").Value = 15878231
$ws.Range("I104").Value = 500
$ws.Range("J104").Value = 16672117
$ws.Range("K104").Value = 1500
$ws.Range("L104").Value = 50016351
$ws.Range("M104").Value = 1121
$ws.Range("N104").Value = -50021593
$ws.Range("H131").Value = 892.5349
$ws.Range("I131").Value = 345.75
$ws.Range("J131").Value = 948.61536
$ws.Range("K131").Value = 1037.25
$ws.Range("L131").Value = 2845.84608
$ws.Range("M131").Value = 4002.75
$ws.Range("N131").Value = -12925.84608

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4105.7144
$ws.Range("I43").Value = 2156
$ws.Range("J43").Value = 8980
$ws.Range("K43").Value = 2156
$ws.Range("L43").Value = 8980
$ws.Range("M43").Value = -2005
$ws.Range("N43").Value = -9282
$ws.Range("H62").Value = 24412.5
$ws.Range("J62").Value = 24412.5
$ws.Range("L62").Value = 24412.5
$ws.Range("N62").Value = -25784.5
$ws.Range("H65").Value = 24412.5
$ws.Range("J65").Value = 24412.5
$ws.Range("L65").Value = 73237.5
$ws.Range("N65").Value = -80101.5
$ws.Range("H92").Value = 875.25
$ws.Range("J92").Value = 875.25
$ws.Range("L92").Value = 875.25
$ws.Range("N92").Value = -4619.25
$ws.Range("H128").Value = 51692
$ws.Range("J128").Value = 51692
$ws.Range("L128").Value = 51692
$ws.Range("N128").Value = -61652

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 9852
$ws.Range("J96").Value = 9852
$ws.Range("L96").Value = 9852
$ws.Range("N96").Value = -15344
$ws.Range("H128").Value = 51698.75
$ws.Range("J128").Value = 51698.75
$ws.Range("L128").Value = 51698.75
$ws.Range("N128").Value = -61658.75
$ws.Range("H133").Value = 52956.7
$ws.Range("J133").Value = 52956.7
$ws.Range("L133").Value = 52956.7
$ws.Range("N133").Value = -58016.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 27998
$ws.Range("J64").Value = 27998
$ws.Range("L64").Value = 27998
$ws.Range("N64").Value = -28494
$ws.Range("H67").Value = 27998
$ws.Range("J67").Value = 27998
$ws.Range("L67").Value = 27998
$ws.Range("N67").Value = -29714
$ws.Range("H93").Value = 26877.8
$ws.Range("J93").Value = 26877.8
$ws.Range("L93").Value = 26877.8
$ws.Range("N93").Value = -31869.8
$ws.Range("H123").Value = 31514.625
$ws.Range("J123").Value = 31514.625
$ws.Range("L123").Value = 31514.625
$ws.Range("N123").Value = -41314.625
$ws.Range("H125").Value = 49878
$ws.Range("J125").Value = 49878
$ws.Range("L125").Value = 49878
$ws.Range("N125").Value = -59718
$ws.Range("H127").Value = 52553.375
$ws.Range("J127").Value = 52553.375
$ws.Range("L127").Value = 52553.375
$ws.Range("N127").Value = -62473.375
$ws.Range("H130").Value = 34740
$ws.Range("J130").Value = 34740
$ws.Range("L130").Value = 34740
$ws.Range("N130").Value = -44780
$ws.Range("H133").Value = 37920.715
$ws.Range("J133").Value = 37920.715
$ws.Range("L133").Value = 37920.715
$ws.Range("N133").Value = -48040.715
